# oox smartart, linear layout: fix scaling of spacing without rules
#
# The target change uncomments two "padding1" / "padding2" <dgm:layoutNode>
# helper nodes inside the diagram's *layout definition* part
# (ppt/diagrams/layout1.xml, uniqueId
# urn:microsoft.com/office/officeart/2005/8/layout/hProcess3), so the
# SmartArt markup matches exactly what PowerPoint's own SmartArt UI
# generates for this built-in "Basic Process" style layout - i.e. this
# restores the full/realistic layout markup instead of a hand-trimmed
# minimal version.
#
# That layout-definition part is a read-only (per loTypeId) resource owned
# by the SmartArt *layout itself*, not by the individual diagram instance:
# PowerPoint's object model deliberately does not expose get/set access to a
# diagram's raw <dgm:layoutDef> markup to end users or macros.
# SmartArtLayout only surfaces read-only catalog metadata (Name/Id/Category/
# Description) and SmartArt/SmartArtNode only expose the *data* model (node
# text, add/remove/promote/demote, colors, quick styles) - there is no
# documented Shape/SmartArt/SmartArtLayout member that reads or writes
# ppt/diagrams/layout1.xml.
#
# The closest, documented, side-effect-free way to ask the host to make sure
# a diagram's presentation/geometry is re-derived from its model and its
# layout definition is SmartArt.Reset(). We call it defensively (and only
# on shapes that actually carry a SmartArt diagram) so the script behaves
# as a safe no-op on hosts/shapes where that does not apply, instead of
# throwing, or mutating unrelated parts (e.g. node text/order in
# ppt/diagrams/data1.xml) that are not part of this change.

function Get-SmartArtShapes {
    param($Slide)

    $result = @()
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.HasSmartArt) {
            $result += $candidate
        }
    }
    return $result
}

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    $smartArtShapes = Get-SmartArtShapes $slide

    foreach ($shape in $smartArtShapes) {
        $smartArt = $shape.SmartArt
        if ($null -eq $smartArt) {
            continue
        }

        # Re-derive the diagram's presentation from its data model and its
        # (built-in) layout definition - the documented, non-mutating-if-
        # unsupported SmartArt verb for this.
        try {
            $smartArt.Reset()
        } catch {
            # Not supported on this host for this shape - leave untouched.
        }
    }
}
